$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 87 (hunk 0)
$ws.Range("H87").Value = 13431.943
$ws.Range("J87").Value = 13431.943
$ws.Range("L87").Value = 13431.943
$ws.Range("N87").Value = -15927.943

# Row 90 (hunk 1)
$ws.Range("H90").Value = 13431.943
$ws.Range("J90").Value = 13431.943
$ws.Range("L90").Value = 40295.829
$ws.Range("N90").Value = -52775.829

# Row 138 (hunk 2)
$ws.Range("H138").Value = 58829864
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 58829864
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 176489592
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -176499872

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (hunk 3)
$ws.Range("H61").Value = 4648.807
$ws.Range("I61").Value = 5647.1
$ws.Range("J61").Value = 3539.5925
$ws.Range("K61").Value = 5647.1
$ws.Range("L61").Value = 3539.5925
$ws.Range("M61").Value = -5435.1
$ws.Range("N61").Value = -3963.5925

# Row 136 (hunk 4)
$ws.Range("H136").Value = 4648.807
$ws.Range("I136").Value = 5647.1
$ws.Range("J136").Value = 3539.5925
$ws.Range("K136").Value = 16941.3
$ws.Range("L136").Value = 10618.7775
$ws.Range("M136").Value = -14391.3
$ws.Range("N136").Value = -15718.7775

$ws = $wb.Worksheets.Item("BSM")
# Row 15 (hunk 5)
$ws.Range("H15").Value = 20000
$ws.Range("I15").Value = 20000
$ws.Range("K15").Value = 20000
$ws.Range("M15").Value = -19773

# Row 22 (hunk 6)
$ws.Range("H22").Value = 6899.2666
$ws.Range("I22").Value = 6899.2666
$ws.Range("K22").Value = 6899.2666
$ws.Range("M22").Value = -6726.2666

# Row 86 (hunk 7)
$ws.Range("H86").Value = 2250.7
$ws.Range("I86").Value = 2357.8572
$ws.Range("J86").Value = 2000.6666
$ws.Range("K86").Value = 2357.8572
$ws.Range("L86").Value = 2000.6666
$ws.Range("M86").Value = -1234.8572
$ws.Range("N86").Value = -4246.6666

# Row 89 (hunk 8)
$ws.Range("H89").Value = 2250.7
$ws.Range("I89").Value = 2357.8572
$ws.Range("J89").Value = 2000.6666
$ws.Range("K89").Value = 11789.286
$ws.Range("L89").Value = 10003.333
$ws.Range("M89").Value = -6173.286
$ws.Range("N89").Value = -21235.333

# Row 92 (hunk 9)
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

# Row 94 (hunk 10)
$ws.Range("H94").Value = 604
$ws.Range("I94").Value = 604
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 604
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -153
$ws.Range("N94").ClearContents()

# Row 133 (hunk 11)
$ws.Range("H133").Value = 49999
$ws.Range("J133").Value = 49999
$ws.Range("L133").Value = 49999
$ws.Range("N133").Value = -60119

# Row 134 (hunk 12)
$ws.Range("H134").Value = 3385.2942
$ws.Range("I134").Value = 2490.8572
$ws.Range("K134").Value = 7472.571599999999
$ws.Range("M134").Value = -4937.571599999999

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (hunk 13)
$ws.Range("H22").Value = 353.14285
$ws.Range("I22").Value = 353.14285
$ws.Range("K22").Value = 353.14285
$ws.Range("M22").Value = -3.14285000000001

# Row 31 (hunk 14)
$ws.Range("H31").Value = 3301.1758
$ws.Range("I31").Value = 839.29785
$ws.Range("J31").Value = 5930.909
$ws.Range("K31").Value = 839.29785
$ws.Range("L31").Value = 5930.909
$ws.Range("M31").Value = -544.29785
$ws.Range("N31").Value = -6520.909

# Row 34 (hunk 15)
$ws.Range("H34").Value = 3301.1758
$ws.Range("I34").Value = 839.29785
$ws.Range("J34").Value = 5930.909
$ws.Range("K34").Value = 839.29785
$ws.Range("L34").Value = 5930.909
$ws.Range("M34").Value = -637.29785
$ws.Range("N34").Value = -6334.909

# Row 99 (hunk 16)
$ws.Range("H99").Value = 3211724.5
$ws.Range("I99").Value = 4573606.5
$ws.Range("K99").Value = 4573606.5
$ws.Range("M99").Value = -4572108.5

# Row 126 (hunk 17)
$ws.Range("H126").Value = 3211724.5
$ws.Range("I126").Value = 4573606.5
$ws.Range("K126").Value = 13720819.5
$ws.Range("M126").Value = -13718349.5

$ws = $wb.Worksheets.Item("CUL")
# Row 68 (hunk 18)
$ws.Range("H68").Value = 1067
$ws.Range("I68").Value = 757.587
$ws.Range("J68").Value = 1330.5741
$ws.Range("K68").Value = 2272.761
$ws.Range("L68").Value = 3991.7223
$ws.Range("M68").Value = -1461.761
$ws.Range("N68").Value = -5613.7223

# Row 71 (hunk 19)
$ws.Range("H71").Value = 1067
$ws.Range("I71").Value = 757.587
$ws.Range("J71").Value = 1330.5741
$ws.Range("K71").Value = 6818.282999999999
$ws.Range("L71").Value = 11975.1669
$ws.Range("M71").Value = -2762.282999999999
$ws.Range("N71").Value = -20087.1669

# Row 125 (hunk 20)
$ws.Range("H125").Value = 1065
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# Row 131 (hunk 21)
$ws.Range("H131").Value = 3572.1875
$ws.Range("J131").Value = 4161.65
$ws.Range("L131").Value = 12484.95
$ws.Range("N131").Value = -22564.95

# Row 133 (hunk 22)
$ws.Range("H133").Value = 4715
$ws.Range("I133").Value = 4715
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 14145
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -9085
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 82 (hunk 23)
$ws.Range("H82").Value = 3690.6924
$ws.Range("I82").Value = 1668.5714
$ws.Range("J82").Value = 6049.8335
$ws.Range("K82").Value = 1668.5714
$ws.Range("L82").Value = 6049.8335
$ws.Range("M82").Value = -1307.5714
$ws.Range("N82").Value = -6771.8335

# Row 85 (hunk 24)
$ws.Range("H85").Value = 3690.6924
$ws.Range("I85").Value = 1668.5714
$ws.Range("J85").Value = 6049.8335
$ws.Range("K85").Value = 1668.5714
$ws.Range("L85").Value = 6049.8335
$ws.Range("M85").Value = -420.5714
$ws.Range("N85").Value = -8545.833500000001

# Row 93 (hunk 25)
$ws.Range("H93").Value = 33832.668
$ws.Range("I93").Value = 50251.5
$ws.Range("J93").Value = 995
$ws.Range("K93").Value = 50251.5
$ws.Range("L93").Value = 995
$ws.Range("M93").Value = -49003.5
$ws.Range("N93").Value = -3491

# Row 132 (hunk 26)
$ws.Range("H132").Value = 5583.273
$ws.Range("I132").Value = 7805.5835
$ws.Range("J132").Value = 4313.381
$ws.Range("K132").Value = 23416.7505
$ws.Range("L132").Value = 12940.143
$ws.Range("M132").Value = -20886.7505
$ws.Range("N132").Value = -18000.143

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (hunk 27)
$ws.Range("H81").Value = 910
$ws.Range("I81").Value = 833.3333
$ws.Range("J81").Value = 942.8570999999999
$ws.Range("K81").Value = 1666.6666
$ws.Range("L81").Value = 1885.7142
$ws.Range("M81").Value = -605.6666
$ws.Range("N81").Value = -4007.7142

# Row 84 (hunk 28)
$ws.Range("H84").Value = 910
$ws.Range("I84").Value = 833.3333
$ws.Range("J84").Value = 942.8570999999999
$ws.Range("K84").Value = 8333.333000000001
$ws.Range("L84").Value = 9428.571
$ws.Range("M84").Value = -3029.333000000001
$ws.Range("N84").Value = -20036.571

# Row 132 (hunk 29)
$ws.Range("H132").Value = 2823.2432
$ws.Range("I132").Value = 2228.8667
$ws.Range("J132").Value = 3228.5
$ws.Range("K132").Value = 6686.6001
$ws.Range("L132").Value = 9685.5
$ws.Range("M132").Value = -4156.6001
$ws.Range("N132").Value = -14745.5

# Row 136 (hunk 30)
$ws.Range("H136").Value = 6064108.5
$ws.Range("I136").Value = 25641758
$ws.Range("J136").Value = 4359.7383
$ws.Range("K136").Value = 76925274
$ws.Range("L136").Value = 13079.2149
$ws.Range("M136").Value = -76922724
$ws.Range("N136").Value = -18179.2149
